# Mifos -> Finflux first-pass changes
# 1) On "Repayment Schedule" sheet, insert a new blank column before column N
#    (this shifts old N -> O, old O -> P, old P -> Q), producing the N/O/P/Q layout seen in the diff.
# 2) Update the active selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

$ws.Columns("N:N").Insert()

$ws.Range("R4").Select()
